$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the anchor paragraph ("Funcionalidad de que guarde el estado del
# switch al cambiar de actividad.") so the two following (empty) paragraphs
# and the later "Semana 3" section can be found reliably, regardless of
# exact paragraph numbering.
# ---------------------------------------------------------------------------

$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Funcionalidad de que guarde el estado del switch al cambiar de actividad.*") {
        $anchorIndex = $i
    }
}

$weekHeadingNewIndex = $anchorIndex + 1   # first of the two empty paragraphs
$bulletNewIndex      = $anchorIndex + 2   # second of the two empty paragraphs
$weekHeadingOldIndex = $anchorIndex + 3   # pre-existing "Semana 3" heading
$bulletOldIndex      = $anchorIndex + 4   # pre-existing GPS/acelerometro bullet

# ---------------------------------------------------------------------------
# 1) The two empty paragraphs right after the anchor become:
#      - a new "Semana 3" bold heading paragraph
#      - a new bulleted list paragraph about GPS / accelerometer data
# ---------------------------------------------------------------------------

$pWeekHeadingNew = $d.Paragraphs.Item($weekHeadingNewIndex)
$rngWeekHeadingNew = $pWeekHeadingNew.Range
$xmlWeekHeadingNew = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>Semana 3</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$rngWeekHeadingNew.InsertXML($xmlWeekHeadingNew)

$pBulletNew = $d.Paragraphs.Item($bulletNewIndex)
$rngBulletNew = $pBulletNew.Range
$xmlBulletNew = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="4"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Añadida la lectura de datos del GPS y del acelerómetro.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$rngBulletNew.InsertXML($xmlBulletNew)

# ---------------------------------------------------------------------------
# 2) The pre-existing "Semana 3" heading further below becomes "Semana 4"
#    (only the digit run's text changes, from "3" to "4"; the "Semana " run
#    must stay untouched/separate, so the single digit character is replaced
#    via a precisely-scoped range rather than a paragraph-wide text/Find
#    replace, which would otherwise merge the two identically-formatted runs
#    into one).
# ---------------------------------------------------------------------------

$pWeekHeadingOld = $d.Paragraphs.Item($weekHeadingOldIndex)
$rngWeekHeadingOld = $pWeekHeadingOld.Range
$digitRange = $d.Range($rngWeekHeadingOld.End - 2, $rngWeekHeadingOld.End - 1)
$xmlDigit = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>4</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$digitRange.InsertXML($xmlDigit)

# ---------------------------------------------------------------------------
# 3) The pre-existing GPS/accelerometer bullet further below gets its text
#    replaced with the new progress note about the wristband / pulse rate.
# ---------------------------------------------------------------------------

$pBulletOld = $d.Paragraphs.Item($bulletOldIndex)
$rngBulletOld = $pBulletOld.Range
$rngBulletOld.Find.Execute( `
    "Añadida la lectura de datos del GPS y del acelerómetro.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Progreso en la implementación de vincular la pulsera con la app, faltando la parte en la que muestre las pulsaciones y reorganizar las actividades.", `
    2)
